$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.778.02"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "'3.525.71"
$ws.Range("E3").Value = "  -3.43%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'196.40"
$ws.Range("E5").Value = "  -3.00%  "
$ws.Range("D6").Value = "'554.53"
$ws.Range("E6").Value = "  -3.15%  "
$ws.Range("D7").Value = "'0.648"
$ws.Range("E7").Value = "  +4.51%  "
$ws.Range("D8").Value = "'3.518.22"
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").Value = "'60.54"
$ws.Range("E11").Value = "  +4.21%  "
$ws.Range("E12").Value = "  -7.10%  "
$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = "  -8.92%  "
$ws.Range("D14").Value = "'9.90"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "'4.091.84"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("D16").Value = "'3.538.04"
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "'67.577.46"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("E20").Value = "  -4.76%  "
$ws.Range("E21").Value = "  -5.42%  "
$ws.Range("D22").Value = "'402.56"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "'87.29"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("E24").Value = "  -5.77%  "
$ws.Range("D25").Value = "'11.71"
$ws.Range("E25").Value = "  -11.44%  "
$ws.Range("D26").Value = "'12.37"
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").Value = "'2.84"
$ws.Range("E27").Value = "  -4.54%  "
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("D30").Value = "'717.75"
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("D31").Value = "'31.39"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "'7.02"
$ws.Range("E32").Value = "  -14.41%  "
$ws.Range("D33").Value = "'11.74"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("D34").Value = "'64.32"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("D36").Value = "'38.48"
$ws.Range("E36").Value = "  -10.31%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  -8.20%  "
$ws.Range("E39").Value = "  -5.09%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "'2.99"
$ws.Range("E41").Value = "  -4.85%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'3.053.83"
$ws.Range("E42").Value = "  -6.91%  "
$ws.Range("D43").Value = "'0.0₃0682"
$ws.Range("E43").Value = "  -13.26%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.49"
$ws.Range("E45").Value = "  -11.84%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.135"
$ws.Range("E46").Value = "  +2.19%  "
$ws.Range("D47").Value = "'0.0408"
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("D48").Value = "'139.64"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "'3.02"
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("E50").Value = "  -16.85%  "
$ws.Range("D51").Value = "'8.25"
$ws.Range("E51").Value = "  -8.09%  "
